# ----------------------------------------------------------------------------
# Edit script for Jogos_da_Semana_FlashScore_2024-11-07.xlsx
#
# Summary of the change:
#   1. Row 2 (ITALY - SERIE A / Genoa-Como): O2 1.44 -> 1.4, P2 2.75 -> 3
#   2. A brand-new match (ARGENTINA - TORNEO BETANO, Defensa y Justicia vs
#      Argentinos Jrs) is inserted as the new row 4, pushing the previously
#      existing rows 4-7 down to rows 5-8 (dimension grows from BD7 to BD8).
#   3. On top of the shift, two of the pushed-down matches also received
#      independent odds updates:
#        - EGYPT - PREMIER LEAGUE (row 5 -> row 6): many odds columns changed
#        - POLAND - DIVISION 1   (row 7 -> row 8): O/P/Q/R odds changed
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update a couple of odds on row 2 ---
$ws.Range("O2").Value = 1.4
$ws.Range("P2").Value = 3

# --- 2) Insert a new row at position 4; existing rows 4-7 shift to 5-8 ---
$ws.Rows.Item(4).Insert()

# --- 3) Fill in the new row 4 (ARGENTINA - TORNEO BETANO match) ---
$ws.Range("A4").Value = "jkXVhfgi"
$ws.Range("B4").Value = "'07/11/2024"
$ws.Range("C4").Value = "19:00"
$ws.Range("D4").Value = "ARGENTINA - TORNEO BETANO"
$ws.Range("E4").Value = "Defensa y Justicia"
$ws.Range("F4").Value = "Argentinos Jrs"
$ws.Range("G4").Value = 2.8
$ws.Range("H4").Value = 2.88
$ws.Range("I4").Value = 2.75
$ws.Range("J4").Value = 3.5
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 3.4
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 8
$ws.Range("O4").Value = 1.4
$ws.Range("P4").Value = 2.75
$ws.Range("Q4").Value = 2.35
$ws.Range("R4").Value = 1.57
$ws.Range("S4").Value = 1.5
$ws.Range("T4").Value = 2.5
$ws.Range("U4").Value = 1.95
$ws.Range("V4").Value = 1.8
$ws.Range("W4").Value = 8
$ws.Range("X4").Value = 13
$ws.Range("Y4").Value = 11
$ws.Range("Z4").Value = 29
$ws.Range("AA4").Value = 26
$ws.Range("AB4").Value = 41
$ws.Range("AC4").Value = 7
$ws.Range("AD4").Value = 5.5
$ws.Range("AE4").Value = 15
$ws.Range("AF4").Value = 51
$ws.Range("AG4").Value = 351
$ws.Range("AH4").Value = 7.5
$ws.Range("AI4").Value = 12
$ws.Range("AJ4").Value = 11
$ws.Range("AK4").Value = 29
$ws.Range("AL4").Value = 23
$ws.Range("AM4").Value = 41
$ws.Range("AN4").Value = 4.5
$ws.Range("AO4").Value = 17
$ws.Range("AP4").Value = 29
$ws.Range("AQ4").Value = 51
$ws.Range("AR4").Value = 81
$ws.Range("AS4").Value = 251
$ws.Range("AT4").Value = 2.5
$ws.Range("AU4").Value = 8.5
$ws.Range("AV4").Value = 67
$ws.Range("AW4").Value = 4.5
$ws.Range("AX4").Value = 15
$ws.Range("AY4").Value = 29
$ws.Range("AZ4").Value = 51
$ws.Range("BA4").Value = 81
$ws.Range("BB4").Value = 251
$ws.Range("BC4").Value = 126
$ws.Range("BD4").Value = 151

# --- 4) Apply the additional odds updates to the EGYPT match, now row 6 ---
$ws.Range("G6").Value = 6.2
$ws.Range("H6").Value = 3.7
$ws.Range("J6").Value = 6
$ws.Range("K6").Value = 2.18
$ws.Range("N6").Value = 6.8
$ws.Range("O6").Value = 1.34
$ws.Range("P6").Value = 3
$ws.Range("Q6").Value = 2
$ws.Range("R6").Value = 1.75
$ws.Range("S6").Value = 1.4
$ws.Range("T6").Value = 2.72
$ws.Range("U6").Value = 2.05
$ws.Range("V6").Value = 1.7
$ws.Range("X6").Value = 37
$ws.Range("Z6").Value = 150
$ws.Range("AA6").Value = 75
$ws.Range("AC6").Value = 6.8
$ws.Range("AD6").Value = 7.2
$ws.Range("AE6").Value = 19.5
$ws.Range("AF6").Value = 110
$ws.Range("AH6").Value = 5.8
$ws.Range("AI6").Value = 6.6
$ws.Range("AJ6").Value = 8.25
$ws.Range("AL6").Value = 13.5
$ws.Range("AT6").Value = 2.72
$ws.Range("AX6").Value = 7.3
$ws.Range("AY6").Value = 17.5
$ws.Range("AZ6").Value = 23

# --- 5) Apply the additional odds updates to the POLAND match, now row 8 ---
$ws.Range("O8").Value = 1.22
$ws.Range("P8").Value = 4
$ws.Range("Q8").Value = 1.75
$ws.Range("R8").Value = 2.05
